$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '51.853.45'
$ws.Range("E2").Value = '  -1.57%  '

# Row 3
$ws.Range("D3").Value = '2.794.21'
$ws.Range("E3").Value = '  -2.00%  '

# Row 4
$ws.Range("E4").Value = '  +0.07%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '358.34'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -0.91%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '109.61'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -6.34%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.558'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +0.99%  '

# Row 8
$ws.Range("E8").Value = '  +0.03%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.591'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -2.56%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '40.01'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -7.05%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0848'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -2.25%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.131'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +0.19%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '19.54'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -3.08%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.65'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -3.32%  '

# Row 15
$ws.Range("D15").Value = '3.238.12'
$ws.Range("E15").Value = '  -1.74%  '

# Row 16
$ws.Range("D16").Value = '2.856.83'
$ws.Range("E16").Value = '  +0.22%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.908'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +0.22%  '

# Row 18
$ws.Range("D18").Value = '51.702.01'
$ws.Range("E18").Value = '  -1.92%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.41'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +1.35%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '3.10'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -3.13%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.10'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -5.12%  '

# Row 22
$ws.Range("D22").Value = '0.0₃0980'
$ws.Range("E22").Value = '  -0.93%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '271.07'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -0.84%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '69.65'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -1.65%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.76'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -2.75%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '26.50'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -3.17%  '

# Row 27
$ws.Range("E27").Value = '  -0.13%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.14'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -2.35%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.23'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -1.68%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.140'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -1.24%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0468'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +2.53%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '51.73'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +1.01%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '33.80'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -2.94%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.75'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -1.89%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.43'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +9.23%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0836'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -0.77%  '

# Row 37
$ws.Range("E37").Value = '  +0.09%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.17'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -3.63%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.00'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -6.47%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '17.76'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -6.19%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.115'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -1.14%  '

# Row 42
$ws.Range("B42").Value = 'Monero'
$ws.Range("C42").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '125.27'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -1.01%  '

# Row 43
$ws.Range("B43").Value = 'Stacks'
$ws.Range("C43").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.50'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -5.70%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.26'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -1.39%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '22.00'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -8.40%  '

# Row 46
$ws.Range("D46").Value = '2.052.54'
$ws.Range("E46").Value = '  -1.51%  '

# Row 47
$ws.Range("B47").Value = 'ApeXProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.32'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +0.35%  '

# Row 48
$ws.Range("B48").Value = 'NEARProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.23'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -4.79%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '5.76'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +1.48%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.930'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -3.95%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '8.99'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -0.44%  '
